# Auto-generated: applies the cryptos.xlsx price/volume/coin update described
# by the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.220.03'
$ws.Range("E2").Value = '  +1.80%  '
$ws.Range("D3").Value = '2.059.98'
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''231.97'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").Value = '''0.622'
$ws.Range("E6").Value = '  +3.36%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''57.02'
$ws.Range("E8").Value = '  +3.61%  '
$ws.Range("D9").Value = '''0.384'
$ws.Range("E9").Value = '  +3.40%  '
$ws.Range("D10").Value = '''57.85'
$ws.Range("E10").Value = '  +0.54%  '
$ws.Range("D11").Value = '''0.0759'
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("E12").Value = '  +1.54%  '
$ws.Range("D13").Value = '2.364.58'
$ws.Range("E13").Value = '  +1.30%  '
$ws.Range("D14").Value = '''14.49'
$ws.Range("E14").Value = '  +1.71%  '
$ws.Range("D15").Value = '''20.77'
$ws.Range("E15").Value = '  +3.70%  '
$ws.Range("D17").Value = '''5.16'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").Value = '2.057.35'
$ws.Range("E18").Value = '  -0.33%  '
$ws.Range("D19").Value = '37.180.44'
$ws.Range("E19").Value = '  +1.21%  '
$ws.Range("D20").Value = '''6.35'
$ws.Range("E20").Value = '  +8.68%  '
$ws.Range("D21").Value = '''69.15'
$ws.Range("E21").Value = '  +1.94%  '
$ws.Range("D22").Value = '0.0₃0808'
$ws.Range("E22").Value = '  +1.44%  '
$ws.Range("D23").Value = '''225.89'
$ws.Range("E23").Value = '  +2.32%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = '''2.42'
$ws.Range("E25").Value = '  +0.84%  '
$ws.Range("D26").Value = '''2.38'
$ws.Range("E26").Value = '  +0.58%  '
$ws.Range("D27").Value = '''165.65'
$ws.Range("E27").Value = '  +1.85%  '
$ws.Range("E28").Value = '  +8.22%  '
$ws.Range("E29").Value = '  +0.38%  '
$ws.Range("E30").Value = '  +0.60%  '
$ws.Range("D31").Value = '''19.05'
$ws.Range("E31").Value = '  +0.90%  '
$ws.Range("E32").Value = '  +0.87%  '
$ws.Range("D33").Value = '''4.44'
$ws.Range("E33").Value = '  +2.07%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '''4.61'
$ws.Range("E34").Value = '  +7.85%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.0617'
$ws.Range("E35").Value = '  +1.71%  '
$ws.Range("D36").Value = '''2.51'
$ws.Range("E36").Value = '  +1.75%  '
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("E38").Value = '  -0.37%  '
$ws.Range("E39").Value = '  +1.22%  '
$ws.Range("D40").Value = '''5.68'
$ws.Range("E40").Value = '  -1.70%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '''96.45'
$ws.Range("E42").Value = '  +3.03%  '
$ws.Range("D43").Value = '1.465.65'
$ws.Range("E43").Value = '  -0.79%  '
$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").Value = '''4.32'
$ws.Range("E44").Value = '  -1.85%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '''1.17'
$ws.Range("E45").Value = '  +5.34%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '''0.0931'
$ws.Range("E46").Value = '  -1.08%  '
$ws.Range("D47").Value = '''0.0212'
$ws.Range("E47").Value = '  +3.95%  '
$ws.Range("E48").Value = '  +1.41%  '
$ws.Range("D49").Value = '''15.10'
$ws.Range("E49").Value = '  -2.92%  '
$ws.Range("E51").Value = '  +1.94%  '

Write-Output "Applied cryptos list update"
